$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.204.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.71%  '

$ws.Range("D3").Value = '''1.917.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.33%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.55%  '

$ws.Range("D5").Value = '''329.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.89%  '

$ws.Range("D6").Value = '''0.9999'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.49%  '

$ws.Range("D7").Value = '''0.5196'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.32%  '

$ws.Range("D8").Value = '''0.4067'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.82%  '

$ws.Range("D9").Value = '''0.08498'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.87%  '

$ws.Range("D10").Value = '''1.128'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.99%  '

$ws.Range("D11").Value = '''42.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.50%  '

$ws.Range("D12").Value = '''23.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +15.08%  '

$ws.Range("D13").Value = '''6.451'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.13%  '

$ws.Range("D14").Value = '''1.920.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.62%  '

$ws.Range("D15").Value = '''7.402'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.99%  '

$ws.Range("E16").Value = '  -0.54%  '

$ws.Range("D17").Value = '''95.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.22%  '

$ws.Range("E18").Value = '  +1.40%  '

$ws.Range("D19").Value = '''0.06688'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.51%  '

$ws.Range("D20").Value = '''18.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.75%  '

$ws.Range("D21").Value = '''0.9996'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.49%  '

$ws.Range("D22").Value = '''6.016'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.36%  '

$ws.Range("D23").Value = '''30.221.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.64%  '

$ws.Range("D24").Value = '''11.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.51%  '

$ws.Range("D25").Value = '''2.235'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.17%  '

$ws.Range("D26").Value = '''2.139.41'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.64%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''21.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.76%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '''161.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.68%  '

$ws.Range("D29").Value = '''2.421'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("D30").Value = '''128.94'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.84%  '

$ws.Range("D31").Value = '''1.101'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.19%  '

$ws.Range("D32").Value = '''0.1068'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.00%  '

$ws.Range("D33").Value = '''6.021'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.26%  '

$ws.Range("D34").Value = '''3.636'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.15%  '

$ws.Range("D35").Value = '''0.02495'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.19%  '

$ws.Range("D36").Value = '''0.06586'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.34%  '

$ws.Range("D37").Value = '''0.2215'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.37%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '''5.189'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.00%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '''1.229'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.45%  '

$ws.Range("D40").Value = '''11.89'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.66%  '

$ws.Range("D41").Value = '''8.830'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.25%  '

$ws.Range("D42").Value = '''0.6532'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.51%  '

$ws.Range("D43").Value = '''1.240'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").Value = '''0.6154'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.81%  '

$ws.Range("D45").Value = '''13.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.67%  '

$ws.Range("D46").Value = '''3.739'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.79%  '

$ws.Range("D47").Value = '''2.081'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.76%  '

$ws.Range("D48").Value = '''1.247'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.67%  '

$ws.Range("D49").Value = '''124.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.57%  '

$ws.Range("D51").Value = '''79.59'
$ws.Range("D51").Style = "Normal"
